$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drag the formula in I5 down through I14 (fill down), same as using the fill handle
$ws.Range("I6:I14").Formula = "=F6-H6"

# Set the active selection to E13, matching the final state after the drag
$ws.Range("E13").Select()
